$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3967.4285
$ws.Range("I62").Value = 3984.889
$ws.Range("J62").Value = 3936
$ws.Range("K62").Value = 3984.889
$ws.Range("L62").Value = 3936
$ws.Range("M62").Value = -3360.889
$ws.Range("N62").Value = -5184
$ws.Range("H65").Value = 3967.4285
$ws.Range("I65").Value = 3984.889
$ws.Range("J65").Value = 3936
$ws.Range("K65").Value = 19924.445
$ws.Range("L65").Value = 19680
$ws.Range("M65").Value = -16804.445
$ws.Range("N65").Value = -25920
$ws.Range("H106").Value = 2449.5833
$ws.Range("I106").Value = 2517.7273
$ws.Range("K106").Value = 2517.7273
$ws.Range("M106").Value = -1886.7273
$ws.Range("H127").Value = 1366.1666
$ws.Range("I127").Value = 1174.25
$ws.Range("J127").Value = 1750
$ws.Range("K127").Value = 3522.75
$ws.Range("L127").Value = 5250
$ws.Range("M127").Value = 1437.25
$ws.Range("N127").Value = -15170
$ws.Range("H137").Value = 1827.4445
$ws.Range("J137").Value = 2374.6667
$ws.Range("L137").Value = 7124.000100000001
$ws.Range("N137").Value = -12224.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3164.6155
$ws.Range("J88").Value = 3460
$ws.Range("L88").Value = 3460
$ws.Range("N88").Value = -4272
$ws.Range("H91").Value = 3164.6155
$ws.Range("J91").Value = 3460
$ws.Range("L91").Value = 3460
$ws.Range("N91").Value = -6268
$ws.Range("H92").Value = 31000
$ws.Range("J92").Value = 31000
$ws.Range("L92").Value = 31000
$ws.Range("N92").Value = -35992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 30386.846
$ws.Range("I82").Value = 15049.2
$ws.Range("J82").Value = 39972.875
$ws.Range("K82").Value = 15049.2
$ws.Range("L82").Value = 39972.875
$ws.Range("M82").Value = -14666.2
$ws.Range("N82").Value = -40738.875
$ws.Range("H85").Value = 30386.846
$ws.Range("I85").Value = 15049.2
$ws.Range("J85").Value = 39972.875
$ws.Range("K85").Value = 15049.2
$ws.Range("L85").Value = 39972.875
$ws.Range("M85").Value = -13723.2
$ws.Range("N85").Value = -42624.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3760
$ws.Range("I12").Value = 586.75
$ws.Range("K12").Value = 586.75
$ws.Range("M12").Value = -416.75
$ws.Range("H31").Value = 2244.1428
$ws.Range("I31").Value = 1651.5652
$ws.Range("J31").Value = 3379.9167
$ws.Range("K31").Value = 1651.5652
$ws.Range("L31").Value = 3379.9167
$ws.Range("M31").Value = -1356.5652
$ws.Range("N31").Value = -3969.9167
$ws.Range("H34").Value = 2244.1428
$ws.Range("I34").Value = 1651.5652
$ws.Range("J34").Value = 3379.9167
$ws.Range("K34").Value = 1651.5652
$ws.Range("L34").Value = 3379.9167
$ws.Range("M34").Value = -1449.5652
$ws.Range("N34").Value = -3783.9167
$ws.Range("H132").Value = 1694
$ws.Range("I132").Value = 1468.25
$ws.Range("K132").Value = 4404.75
$ws.Range("M132").Value = -1874.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 377.25
$ws.Range("I7").Value = 212
$ws.Range("J7").Value = 542.5
$ws.Range("K7").Value = 636
$ws.Range("L7").Value = 1627.5
$ws.Range("M7").Value = -524
$ws.Range("N7").Value = -1851.5
$ws.Range("H55").Value = 2864.7
$ws.Range("I55").Value = 1283.3334
$ws.Range("J55").Value = 3542.4285
$ws.Range("K55").Value = 3850.0002
$ws.Range("L55").Value = 10627.2855
$ws.Range("M55").Value = -3673.0002
$ws.Range("N55").Value = -10981.2855
$ws.Range("H120").Value = 5033
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H131").Value = 943.75
$ws.Range("I131").Value = 925
$ws.Range("K131").Value = 2775
$ws.Range("M131").Value = 2265
$ws.Range("H137").Value = 907.5
$ws.Range("I137").Value = 907.5
$ws.Range("K137").Value = 2722.5
$ws.Range("M137").Value = 2377.5
$ws.Range("H138").Value = 5468.6
$ws.Range("I138").Value = 4830
$ws.Range("J138").Value = 6426.5
$ws.Range("K138").Value = 14490
$ws.Range("L138").Value = 19279.5
$ws.Range("M138").Value = -9350
$ws.Range("N138").Value = -29559.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3793.6
$ws.Range("I80").Value = 3826.6667
$ws.Range("J80").Value = 3744
$ws.Range("K80").Value = 3826.6667
$ws.Range("L80").Value = 3744
$ws.Range("M80").Value = -2828.6667
$ws.Range("N80").Value = -5740
$ws.Range("H83").Value = 3793.6
$ws.Range("I83").Value = 3826.6667
$ws.Range("J83").Value = 3744
$ws.Range("K83").Value = 19133.3335
$ws.Range("L83").Value = 18720
$ws.Range("M83").Value = -14141.3335
$ws.Range("N83").Value = -28704
$ws.Range("H102").Value = 2594.5
$ws.Range("I102").Value = 2594.5
$ws.Range("K102").Value = 2594.5
$ws.Range("M102").Value = -972.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 3144.7742
$ws.Range("I46").Value = 2280.6875
$ws.Range("J46").Value = 4066.4666
$ws.Range("K46").Value = 2280.6875
$ws.Range("L46").Value = 4066.4666
$ws.Range("M46").Value = -2092.6875
$ws.Range("N46").Value = -4442.4666
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H99").Value = 31749
$ws.Range("I99").Value = 31749
$ws.Range("K99").Value = 31749
$ws.Range("M99").Value = -28754
$ws.Range("H122").Value = 3489.3845
$ws.Range("I122").Value = 3305.6365
$ws.Range("K122").Value = 9916.9095
$ws.Range("M122").Value = -7466.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 2286.1
$ws.Range("J122").Value = 2652.5
$ws.Range("L122").Value = 7957.5
$ws.Range("N122").Value = -12857.5
